$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 0.09420142906393127
$ws.Range("B3").Value = 0.003048912659270201
$ws.Range("C3").Value = 0.0005962938982849801
$ws.Range("D3").Value = 3.710903485058106
$ws.Range("E3").Value = 0.008633298165580937
$ws.Range("F3").Value = 0.001880192533460522
$ws.Range("G3").Value = 0.004217632785079879
$ws.Range("H3").Value = 0.09725034172320146
$ws.Range("B4").Value = 0.006161054489117924
$ws.Range("C4").Value = 0.0009669337875731515
$ws.Range("D4").Value = 6.544345097731695
$ws.Range("E4").Value = 0.1307511425145125
$ws.Range("F4").Value = 0.004265890102318803
$ws.Range("G4").Value = 0.008056218875917048
$ws.Range("H4").Value = 0.1003624835530492
$ws.Range("B5").Value = 0.0139757799171559
$ws.Range("C5").Value = 0.002128987993842845
$ws.Range("D5").Value = 10.27903694035477
$ws.Range("E5").Value = 0.03011572686101117
$ws.Range("F5").Value = 0.009803019988957475
$ws.Range("G5").Value = 0.01814853984535432
$ws.Range("H5").Value = 0.1081772089810872
$ws.Range("B6").Value = 0.02769360922876226
$ws.Range("C6").Value = 0.003859606362895035
$ws.Range("D6").Value = 9.463451900014437
$ws.Range("E6").Value = 0.05392725218646769
$ws.Range("F6").Value = 0.02012889306228829
$ws.Range("G6").Value = 0.03525832539523622
$ws.Range("H6").Value = 0.1218950382926935
$ws.Range("B7").Value = 0.05823091844794599
$ws.Range("C7").Value = 0.008240921759839624
$ws.Range("D7").Value = 11.74256019774642
$ws.Range("E7").Value = 0.04604176990406618
$ws.Range("F7").Value = 0.04207896355436509
$ws.Range("G7").Value = 0.07438287334152691
$ws.Range("H7").Value = 0.1524323475118773
$ws.Range("B8").Value = 0.04841200587295257
$ws.Range("C8").Value = 0.005392733206263865
$ws.Range("D8").Value = 14.6370017898739
$ws.Range("E8").Value = 0.05505200337544913
$ws.Range("F8").Value = 0.03784240638950259
$ws.Range("G8").Value = 0.05898160535640255
$ws.Range("H8").Value = 0.1426134349368838
$ws.Range("B9").Value = 0.0635275368151086
$ws.Range("C9").Value = 0.00719214180242916
$ws.Range("D9").Value = 12.0268768325909
$ws.Range("E9").Value = 0.1208583929595987
$ws.Range("F9").Value = 0.04943115598679013
$ws.Range("G9").Value = 0.07762391764342709
$ws.Range("H9").Value = 0.1577289658790399
$ws.Range("B10").Value = -0.09420142906393127
$ws.Range("C10").Value = 0.0004998572668878386
$ws.Range("D10").Value = -226.7623461832002
$ws.Range("E10").Value = 0
$ws.Range("F10").Value = -0.09518113603888863
$ws.Range("G10").Value = -0.09322172208897389
$ws.Range("B11").Value = -0.04279039862578433
$ws.Range("C11").Value = 0.0005550383961205047
$ws.Range("D11").Value = -89.73190822408711
$ws.Range("E11").Value = 0.000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000006856851286177836
$ws.Range("F11").Value = -0.04387825921392318
$ws.Range("G11").Value = -0.04170253803764547
$ws.Range("H11").Value = 0.05141103043814694
$ws.Range("B12").Value = -0.03554855962313918
$ws.Range("C12").Value = 0.000543600921371943
$ws.Range("D12").Value = -76.39585559168611
$ws.Range("E12").Value = 0.0000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000004387874846556969
$ws.Range("F12").Value = -0.03661400306549098
$ws.Range("G12").Value = -0.03448311618078738
$ws.Range("H12").Value = 0.05865286944079209
$ws.Range("B13").Value = -0.03186218027827979
$ws.Range("C13").Value = 0.0005369305992643614
$ws.Range("D13").Value = -67.33658055144646
$ws.Range("E13").Value = 0.0000000000000004184231528168113
$ws.Range("F13").Value = -0.03291455001765726
$ws.Range("G13").Value = -0.03080981053890231
$ws.Range("H13").Value = 0.06233924878565147
$ws.Range("B14").Value = -0.02830834178911806
$ws.Range("C14").Value = 0.000531012311787611
$ws.Range("D14").Value = -59.06220349157483
$ws.Range("E14").Value = 0.00000000002142622865227341
$ws.Range("F14").Value = -0.0293491116975341
$ws.Range("G14").Value = -0.02726757188070202
$ws.Range("H14").Value = 0.06589308727481322
$ws.Range("B15").Value = -0.02214728938908753
$ws.Range("C15").Value = 0.000518877148626403
$ws.Range("D15").Value = -49.35781968837958
$ws.Range("E15").Value = 0.00003867555779201483
$ws.Range("F15").Value = -0.0231642747461753
$ws.Range("G15").Value = -0.02113030403199976
$ws.Range("H15").Value = 0.07205413967484373
$ws.Range("B16").Value = -0.02206827084049732
$ws.Range("C16").Value = 0.0005092709938389442
$ws.Range("D16").Value = -49.45060715153529
$ws.Range("E16").Value = 0.0160739469037296
$ws.Range("F16").Value = -0.02306642852385788
$ws.Range("G16").Value = -0.02107011315713676
$ws.Range("H16").Value = 0.07213315822343394
$ws.Range("B17").Value = -0.01718308746102696
$ws.Range("C17").Value = 0.0005114946841900559
$ws.Range("D17").Value = -38.80609499517521
$ws.Range("E17").Value = 0.0007098777829041644
$ws.Range("F17").Value = -0.01818560338147268
$ws.Range("G17").Value = -0.01618057154058125
$ws.Range("H17").Value = 0.0770183416029043
$ws.Range("B18").Value = -0.01601625246254909
$ws.Range("C18").Value = 0.0005196653275295691
$ws.Range("D18").Value = -35.92931741135929
$ws.Range("E18").Value = 0.000007647475809324413
$ws.Range("F18").Value = -0.01703478265694762
$ws.Range("G18").Value = -0.01499772226815056
$ws.Range("H18").Value = 0.07818517660138218
$ws.Range("B19").Value = -0.01348634223331539
$ws.Range("C19").Value = 0.0005152318884402656
$ws.Range("D19").Value = -30.30745393491721
$ws.Range("E19").Value = 0.02555701692463587
$ws.Range("F19").Value = -0.01449618299262599
$ws.Range("G19").Value = -0.0124765014740048
$ws.Range("H19").Value = 0.08071508683061587
$ws.Range("B20").Value = -0.01215317187024798
$ws.Range("C20").Value = 0.0005369545297122776
$ws.Range("D20").Value = -26.79104244387864
$ws.Range("E20").Value = 0.05007558715002679
$ws.Range("F20").Value = -0.01320558859908848
$ws.Range("G20").Value = -0.01110075514140747
$ws.Range("H20").Value = 0.08204825719368329
$ws.Range("B21").Value = -0.01114775004535714
$ws.Range("C21").Value = 0.0005511354133683927
$ws.Range("D21").Value = -22.56049699506488
$ws.Range("E21").Value = 0.04714242947776345
$ws.Range("F21").Value = -0.0122279609216971
$ws.Range("G21").Value = -0.01006753916901718
$ws.Range("H21").Value = 0.08305367901857413
$ws.Range("B22").Value = -0.009095904215864909
$ws.Range("C22").Value = 0.0005450284230243712
$ws.Range("D22").Value = -18.18755923677869
$ws.Range("E22").Value = 0.05201737307034413
$ws.Range("F22").Value = -0.0101641454803702
$ws.Range("G22").Value = -0.008027662951359617
$ws.Range("H22").Value = 0.08510552484806636
$ws.Range("B23").Value = -0.008115053144098229
$ws.Range("C23").Value = 0.0005511406855162885
$ws.Range("D23").Value = -15.65538627201859
$ws.Range("E23").Value = 0.001218237129286086
$ws.Range("F23").Value = -0.009195274299665946
$ws.Range("G23").Value = -0.007034831988530505
$ws.Range("H23").Value = 0.08608637591983304
$ws.Range("B24").Value = -0.006932632782908826
$ws.Range("C24").Value = 0.0005462834129494593
$ws.Range("D24").Value = -13.14406989324688
$ws.Range("E24").Value = 0.006044530511736924
$ws.Range("F24").Value = -0.00800333380136835
$ws.Range("G24").Value = -0.005861931764449301
$ws.Range("H24").Value = 0.08726879628102244
$ws.Range("B25").Value = -0.004423307554967946
$ws.Range("C25").Value = 0.0005168599487469275
$ws.Range("D25").Value = -9.383190213251124
$ws.Range("E25").Value = 0.0000142093771019859
$ws.Range("F25").Value = -0.00543633926985044
$ws.Range("G25").Value = -0.003410275840085452
$ws.Range("H25").Value = 0.08977812150896332
$ws.Range("B26").Value = 0.160890112647672
$ws.Range("C26").Value = 0.007761679762373025
$ws.Range("D26").Value = 45.39137576020961
$ws.Range("E26").Value = 0.02047360810682636
$ws.Range("F26").Value = 0.145677457287413
$ws.Range("G26").Value = 0.1761027680079308
$ws.Range("H26").Value = 0.2550915417116033
